$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.485.10'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.744.35'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.29%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.25%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4221'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -9.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3580'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.59'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07422'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.113'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.19%  '
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.104'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.184'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.742.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001065'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.59'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06180'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.24%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.101'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5248'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '27.513.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.327'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.362'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.939.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '126.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.202'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.687'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09135'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.691'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.66'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02290'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.087'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.29%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2125'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06081'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6392'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.191'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.422'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.899'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.719'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5871'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.951'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06848'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.42%  '
